$wb = $excel.ActiveWorkbook

# --- Sheet "bar" (sheet1): select A1:C3, but it will stop being the active/tabSelected sheet ---
$bar = $wb.Worksheets.Item("bar")
$bar.Range("A1:C3").Select()

# --- Sheet "foo" (sheet2): fill in new data. ---
# Enter values in an order that reproduces the target shared-string table order:
# karina(7), kayla(8), baby(9), mom(10), name(11), title(12), og(13)
$foo = $wb.Worksheets.Item("foo")

$foo.Range("A2").Value = "karina"
$foo.Range("A3").Value = "kayla"

$foo.Range("B2").Value = "baby"
$foo.Range("B3").Value = "mom"

$foo.Range("A1").Value = "name"
$foo.Range("B1").Value = "title"
$foo.Range("C1").Value = "og"

$foo.Range("C2").Value = 1
$foo.Range("C3").Value = 0

# Make "foo" the active tab (activeTab="1" in workbook.xml, tabSelected on sheet2)
$foo.Activate()
